$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 120, pushing existing row 120 (and below) down to row 121+
$ws.Rows.Item(120).Insert()

# Populate the newly inserted row 120 with the new data record
$ws.Range("A120").Value = 9
$ws.Range("B120").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C120").Value = "Metropolitana"
$ws.Range("D120").Value = Get-Date -Year 2022 -Month 1 -Day 5 -Hour 0 -Minute 0 -Second 0
$ws.Range("E120").Value = 13
$ws.Range("F120").Value = 100112021
$ws.Range("G120").Value = "Ají"
$ws.Range("H120").Value = "Americana (o)"
$ws.Range("I120").Value = "Primera"
$ws.Range("J120").Value = 7
$ws.Range("K120").Value = 25000
$ws.Range("L120").Value = 26000
$ws.Range("M120").Value = 25429
$ws.Range("N120").Value = "$/caja 15 kilos"
$ws.Range("O120").Value = "Provincia de Huasco"
$ws.Range("P120").Value = 1695
$ws.Range("Q120").Value = 15
$ws.Range("R120").Value = "Hortaliza"
